# sua ban ton zalo
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new incident data
$ws.Range("A2").Value = "UL_TTT093M_HNI"
$ws.Range("B2").Value = "THACH-HOA-TTT_HNI"
$ws.Range("C2").Value = "POWER_AC_EAS"
$ws.Range("D2").Value = "07/05/2025 13:18:36"
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "Trạm viễn thông loại 3"
$ws.Range("G2").Value = "Thạch Thất"
$ws.Range("H2").Value = 0.38

# Update row 3 with the new incident data
$ws.Range("A3").Value = "3G_DPG045S_HNI"
$ws.Range("B3").Value = "TRUNG-CHAU-VAN-MON2-11-SMC-DPG_HNI"
$ws.Range("C3").Value = "SITE_OOS"
$ws.Range("D3").Value = "06/05/2025 23:18:18"
$ws.Range("E3").Value = "Trạm smc mất điện - 1 - sonnn - 06/05/2025 23:24:34"
$ws.Range("F3").Value = "Trạm viễn thông loại 3"
$ws.Range("G3").Value = "Đan Phượng"
$ws.Range("H3").Value = 14.39

# Remove the now-obsolete rows 4-8 (old incidents no longer reported)
$ws.Range("A4:H8").EntireRow.Delete()

# Adjust column widths to fit the new, shorter content
$ws.Range("B1").ColumnWidth = 36
$ws.Range("C1").ColumnWidth = 14
$ws.Range("E1").ColumnWidth = 53
